# Insert a new weekly price record as row 11 (Fecha serial 44414 = 2021-08-06),
# pushing the existing rows 11-45 down to 12-46 (dimension grows to A1:R46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("11").Insert()

$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "Femacal de La Calera"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44414
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 100112026
$ws.Range("G11").Value = "Haba"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 40
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 14000
$ws.Range("N11").Value = "`$/saco 25 kilos"
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 560
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = "Hortaliza"
